$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 49; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    if ($cell.Text -eq "OD600") {
        $cell.Value = "Optical Density"
    }
}
